# Apply the "Atualizado por script em 10-11-2023 14:45" update to the
# Bulgaria vtora-liga 2023-2024 odds sheet.
#
# The match-detail columns (F:V) for a handful of rows were reshuffled
# (the underlying scrape re-ordered same-date fixtures) and one brand new
# fixture (CSKA 1948 Sofia II vs Belasitsa, played 10/11/2023) was
# appended as row 143. Columns A:E (Indice / pais / torneio / temporada /
# data_partida) stay put for every row position - only the F:V payload
# moves between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Snapshot the F:V payload of every row that is either a source or a
#    destination of a move, BEFORE any writes happen (some rows are both
#    a source for one move and a destination for another).
# ---------------------------------------------------------------------
$rowsToSnapshot = 126,133,134,135,136,137,138,141,142
$snapshot = @{}
foreach ($r in $rowsToSnapshot) {
    $snapshot[$r] = $ws.Range("F$r`:V$r").Value2
}

# ---------------------------------------------------------------------
# 2) Write each destination row with the snapshotted payload of its
#    source row.
# ---------------------------------------------------------------------
$moves = @{
    126 = 133   # row 126 <- old row 133 payload
    133 = 126   # row 133 <- old row 126 payload
    134 = 135   # row 134 <- old row 135 payload
    135 = 136   # row 135 <- old row 136 payload
    136 = 137   # row 136 <- old row 137 payload
    137 = 138   # row 137 <- old row 138 payload
    138 = 134   # row 138 <- old row 134 payload
    141 = 142   # row 141 <- old row 142 payload
    142 = 141   # row 142 <- old row 141 payload
}

foreach ($dest in $moves.Keys) {
    $src = $moves[$dest]
    $ws.Range("F$dest`:V$dest").Value2 = $snapshot[$src]
}

# ---------------------------------------------------------------------
# 3) Append the brand new fixture as row 143, matching the formatting
#    used by every other data row (bold/bordered/centered index cell in
#    A, date-formatted cell in E, plain cells for the rest).
# ---------------------------------------------------------------------
$newRow = 143

$ws.Range("A126").Copy()
$ws.Range("A$newRow").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E126").Copy()
$ws.Range("E$newRow").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("A$newRow").Value = 142
$ws.Range("B$newRow").Value = "bulgaria"
$ws.Range("C$newRow").Value = "vtora-liga"
$ws.Range("D$newRow").Value = "2023-2024"
$ws.Range("E$newRow").Value = 45240.5625
$ws.Range("F$newRow").Value = "CSKA 1948 Sofia II"
$ws.Range("G$newRow").Value = 0
$ws.Range("H$newRow").Value = "Belasitsa"
$ws.Range("I$newRow").Value = 0
$ws.Range("J$newRow").Value = 1.35
$ws.Range("K$newRow").Value = "10/11/2023 02:42"
$ws.Range("L$newRow").Value = 1.5
$ws.Range("M$newRow").Value = "10/11/2023 13:29"
$ws.Range("N$newRow").Value = 4.28
$ws.Range("O$newRow").Value = "10/11/2023 02:42"
$ws.Range("P$newRow").Value = 3.57
$ws.Range("Q$newRow").Value = "10/11/2023 13:29"
$ws.Range("R$newRow").Value = 7.38
$ws.Range("S$newRow").Value = "10/11/2023 02:42"
$ws.Range("T$newRow").Value = 6.26
$ws.Range("U$newRow").Value = "10/11/2023 13:29"
$ws.Range("V$newRow").Value = "https://www.betexplorer.com/football/bulgaria/vtora-liga/cska-1948-sofia-belasitsa-petrich/pS69BmHg/"
